$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name (shared between both sheets' B1 cell) to the new value
$newName = "4247-RBI-EI-DB-SAR-REC-RNI-FEE+INTEREST-FFConMONonLASTSUN-FIFC-1-FFROP-DL-FIFR-1-MD-TR-1-ONT-PE-1st"
$ws1.Range("B1").Value = $newName
$ws2.Range("B1").Value = $newName

# Change shortname (B2) from numeric 4247 to the text value "424d"
$ws1.Range("B2").Value = "424d"

# Move the selection on the input sheet to B9
$ws1.Range("B9").Select()

# Make the output sheet the active/selected tab
$ws2.Activate()
$ws2.Range("B1").Select()
